$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A94:E94")
$newRow.NumberFormat = "@"

$ws.Range("A94").Value = "2025-12-19"
$ws.Range("B94").Value = "Pick 3"
$ws.Range("C94").Value = "251219"
$ws.Range("D94").Value = "7-7-9"
$ws.Range("E94").Value = "2025-12-19T21:37:22.458+04:00"
